# Processed Results - mean and median increase
# Adds "Mean increase" / "Median increase" computed columns (D/F) to the
# Low / Medium / High frequency summary blocks and an overall summary in
# the combined "All" block, matching the author's manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-IncreaseBlock {
    param($HeaderRow, $FormulaRow, $MeanRef, $MedianRef)

    $dHeader = $ws.Range("D$HeaderRow")
    $dHeader.Value = "Mean increase"
    $dHeader.Font.Bold = $true

    $fHeader = $ws.Range("F$HeaderRow")
    $fHeader.Value = "Median increase"
    $fHeader.Font.Bold = $true

    $dFormula = $ws.Range("D$FormulaRow")
    $dFormula.Formula = "=(($MeanRef / 95.321842) * 100) - 100"
    $dFormula.Style = "Standaard"

    $fFormula = $ws.Range("F$FormulaRow")
    $fFormula.Formula = "=(($MedianRef / 95.22216) * 100) - 100"
    $fFormula.Style = "Standaard"
}

# Low frequency block (data rows 2-31; mean in E3, median in E10)
Add-IncreaseBlock 18 19 "E3" "E10"

# Medium frequency block (data rows 34-63; mean in E35, median in E42)
Add-IncreaseBlock 50 51 "E35" "E42"

# High frequency block (data rows 66-95; mean in E67, median in E74)
Add-IncreaseBlock 82 83 "E67" "E74"

# Overall "All" summary (headers only; formula row averages the three blocks)
$dHeaderAll = $ws.Range("D113")
$dHeaderAll.Value = "Mean increase"
$dHeaderAll.Font.Bold = $true

$fHeaderAll = $ws.Range("F113")
$fHeaderAll.Value = "Median increase"
$fHeaderAll.Font.Bold = $true

$dAll = $ws.Range("D114")
$dAll.Formula = "=(D19 + D51 + D83) / 3"
$dAll.Style = "Standaard"

$fAll = $ws.Range("F114")
$fAll.Formula = "=(F19 + F51 + F83) / 3"
$fAll.Style = "Standaard"

# --- Rebuild the _xlchart.v1.* defined names so duplicate entries are
# dropped and the survivors are renumbered/reordered exactly like Excel
# does when it re-saves the workbook after the chart source ranges change.
$wb.Names.Item("_xlchart.v1.5").Delete()   # duplicate of v1.10 (A98:A187)
$wb.Names.Item("_xlchart.v1.6").Delete()   # duplicate of v1.11 (B98:B187)

$wb.Names.Item("_xlchart.v1.3").Name = "_xlchart.v1.tmp8"
$wb.Names.Item("_xlchart.v1.4").Name = "_xlchart.v1.tmp9"
$wb.Names.Item("_xlchart.v1.7").Name = "_xlchart.v1.tmp3"
$wb.Names.Item("_xlchart.v1.8").Name = "_xlchart.v1.tmp4"
$wb.Names.Item("_xlchart.v1.9").Name = "_xlchart.v1.tmp5"
$wb.Names.Item("_xlchart.v1.10").Name = "_xlchart.v1.tmp6"
$wb.Names.Item("_xlchart.v1.11").Name = "_xlchart.v1.tmp7"

$wb.Names.Item("_xlchart.v1.tmp3").Name = "_xlchart.v1.3"
$wb.Names.Item("_xlchart.v1.tmp4").Name = "_xlchart.v1.4"
$wb.Names.Item("_xlchart.v1.tmp5").Name = "_xlchart.v1.5"
$wb.Names.Item("_xlchart.v1.tmp6").Name = "_xlchart.v1.6"
$wb.Names.Item("_xlchart.v1.tmp7").Name = "_xlchart.v1.7"
$wb.Names.Item("_xlchart.v1.tmp8").Name = "_xlchart.v1.8"
$wb.Names.Item("_xlchart.v1.tmp9").Name = "_xlchart.v1.9"

# Re-order the defined names list so it is written out sorted exactly
# like the target file (v1.0 .. v1.9 in ascending order).
$orderedNames = @(
    "_xlchart.v1.0", "_xlchart.v1.1", "_xlchart.v1.2", "_xlchart.v1.3",
    "_xlchart.v1.4", "_xlchart.v1.5", "_xlchart.v1.6", "_xlchart.v1.7",
    "_xlchart.v1.8", "_xlchart.v1.9"
)
$savedRefs = @{}
foreach ($nm in $orderedNames) {
    $savedRefs[$nm] = $wb.Names.Item($nm).RefersTo
}
foreach ($nm in $orderedNames) {
    $wb.Names.Item($nm).Delete()
}
foreach ($nm in $orderedNames) {
    $added = $wb.Names.Add($nm, $savedRefs[$nm])
    $added.Visible = $false
}

# Move the sheet's selection/scroll position close to where the new data
# was added, matching the state Excel saved the workbook in.
$ws.Range("I114").Select()
